$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text so numeric-looking price strings (e.g. '220.72')
# are stored as text, matching the original inlineStr cell type, not as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = '26.321.28'
$ws.Cells.Item(2, 5).Value = '  +0.58%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.667.02'
$ws.Cells.Item(3, 5).Value = '  +0.83%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '220.72'

# Row 6
$ws.Cells.Item(6, 4).Value = '0.5311'
$ws.Cells.Item(6, 5).Value = '  +0.07%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.38%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '0.2650'
$ws.Cells.Item(8, 5).Value = '  +1.04%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '0.06366'

# Row 10
$ws.Cells.Item(10, 4).Value = '20.84'
$ws.Cells.Item(10, 5).Value = '  +2.15%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '0.07853'
$ws.Cells.Item(11, 5).Value = '  +0.56%  '

# Row 12
$ws.Cells.Item(12, 2).Value = 'Polkadot'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(12, 4).Value = '4.515'
$ws.Cells.Item(12, 5).Value = '  -0.10%  '

# Row 13
$ws.Cells.Item(13, 2).Value = 'WrappedEther'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(13, 4).Value = '1.673.44'
$ws.Cells.Item(13, 5).Value = '  +1.30%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '1.895.57'
$ws.Cells.Item(14, 5).Value = '  +0.80%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '0.5592'
$ws.Cells.Item(15, 5).Value = '  +1.80%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '0.0₅8157'
$ws.Cells.Item(16, 5).Value = '  -0.20%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '65.85'
$ws.Cells.Item(17, 5).Value = '  +0.63%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '26.316.70'
$ws.Cells.Item(18, 5).Value = '  +0.67%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +0.42%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  +2.49%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '197.02'
$ws.Cells.Item(21, 5).Value = '  +3.02%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '10.26'
$ws.Cells.Item(22, 5).Value = '  +1.78%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '6.043'
$ws.Cells.Item(23, 5).Value = '  +0.65%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  +0.34%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '145.43'
$ws.Cells.Item(25, 5).Value = '  +0.18%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '0.1217'
$ws.Cells.Item(26, 5).Value = '  -0.40%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '7.235'
$ws.Cells.Item(27, 5).Value = '  +0.48%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '16.14'
$ws.Cells.Item(28, 5).Value = '  +0.99%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +2.15%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '0.05883'
$ws.Cells.Item(30, 5).Value = '  +2.55%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +0.79%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -0.31%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '3.332'
$ws.Cells.Item(33, 5).Value = '  +2.10%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '1.605'
$ws.Cells.Item(34, 5).Value = '  +0.96%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'MXToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(35, 4).Value = '2.828'
$ws.Cells.Item(35, 5).Value = '  +0.74%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'ARBITRUM'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(36, 4).Value = '0.9608'
$ws.Cells.Item(36, 5).Value = '  +1.25%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '2.436'
$ws.Cells.Item(37, 5).Value = '  +0.57%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '0.5808'
$ws.Cells.Item(38, 5).Value = '  +1.14%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '0.01612'
$ws.Cells.Item(39, 5).Value = '  +0.60%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '5.942'
$ws.Cells.Item(40, 5).Value = '  +2.71%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '1.074.14'
$ws.Cells.Item(41, 5).Value = '  +3.39%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '0.8618'
$ws.Cells.Item(42, 5).Value = '  +1.40%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  +0.42%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '102.76'
$ws.Cells.Item(44, 5).Value = '  -1.07%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '1.806.34'
$ws.Cells.Item(45, 5).Value = '  +0.71%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '58.34'
$ws.Cells.Item(46, 5).Value = '  +2.66%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  +1.33%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '1.016'
$ws.Cells.Item(48, 5).Value = '  +1.02%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '0.4412'
$ws.Cells.Item(49, 5).Value = '  +1.29%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '7.994'
$ws.Cells.Item(50, 5).Value = '  +1.71%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '0.05151'
$ws.Cells.Item(51, 5).Value = '  -0.06%  '

# Restore default (Normal) style on column D so no stray number-format style
# is left attached to the cells (keeps XML identical to original styling).
$ws.Range("D2:D51").Style = "Normal"
